# "Add stage 3 (B1-B2)" - fill in the TxHash values on the B1 and B2
# worksheets (the previously-blank "stage 3" rows) and leave the
# workbook/window selection state the way Excel would after the user
# worked through Info -> B1 -> B2 in that order (B2 ending up active).

$wb = $excel.ActiveWorkbook

# --- Info sheet: visit it first (matches its existing C10 selection),
# this also means it will NOT be the final active/tabSelected sheet.
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Activate()
$wsInfo.Range("C10").Select()

# --- B1: fill in the two new TxHash rows for stage 3 and leave the
# selection on A4 (the next empty row).
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Activate()
$wsB1.Range("A2").Value = "C16D20B44314B3FBEAC198EDDA69A738937AE90AE5F4E0BD64104F321BF82FB3"
$wsB1.Range("A3").Value = "AAADE957FDF44AEF092AC887F3ACAB1EA1E12CFC11E9E9EC483282055FA116CB"
$wsB1.Range("A4").Select()

# --- B2: fill in its two new TxHash rows for stage 3. The second one
# (A3) picked up a stray number-format tweak (scientific notation) in
# the original edit, so reproduce that too. Leave the selection on E7,
# and finish here so B2 stays the active sheet/tab.
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Activate()
$wsB2.Range("A2").Value = "8D377032E6CBEFC276C68FF3AA35C1126474E26D31CABFDE138BB2FF46D32271"
$wsB2.Range("A3").Value = "85E290743B46237632CD731B55E4987732757AA29A8415FE3D610CAF4E92ED0C"
$wsB2.Range("A3").NumberFormat = "0.00E+00"
$wsB2.Range("E7").Select()
